$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.5400586209180323, 1.663049718235778, -1.117886511442075, 2.513879836028467)"
$ws.Range("C2").Value = "NCT(2.875718548035966, 1.9242384096750262, -0.015852927372324348, 4.064348798436132)"
$ws.Range("D2").Value = "JSU(-0.8199789554989656, 1.013170099287664, 0.5008443120330275, 2.6918353655703955)"
$ws.Range("E2").Value = "JSU(-0.7797907705206331, 1.075815679049899, 4.797260706476486, 4.104736785649902)"
